# Insert a new row 17 ("click" / "the first Use this address button")
# which pushes the former row 17 (assert / h2#deliver-to-customer-text /
# Delivering to Normar Weaver) down to row 18, then select C18 to match
# the saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 17, shifting rows 17+ down.
$ws.Rows.Item(17).Insert() | Out-Null

# Populate the newly inserted row 17.
$ws.Range("A17").Value = "click"
$ws.Range("B17").Value = "the first Use this address button"

# The row-insert operation carries the formatted-but-empty E column cell
# down into the new row; remove it so the row only spans A:B like the
# target workbook.
$ws.Range("E17").Clear() | Out-Null

# Match the workbook's saved selection.
$ws.Range("C18").Select() | Out-Null
